$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.860.76"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").Value = "1.906.54"
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.42%  "

$ws.Range("D5").Value = "'313.10"
$ws.Range("E5").Value = "  -0.70%  "

$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.41%  "

$ws.Range("D7").Value = "'0.5001"
$ws.Range("E7").Value = "  +4.09%  "

$ws.Range("D8").Value = "'0.3815"
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "'0.07282"
$ws.Range("E9").Value = "  -1.09%  "

$ws.Range("D10").Value = "'0.9107"
$ws.Range("E10").Value = "  -2.46%  "

$ws.Range("E11").Value = "  +0.77%  "

$ws.Range("D12").Value = "'0.07638"
$ws.Range("E12").Value = "  -1.79%  "

$ws.Range("D13").Value = "1.880.73"
$ws.Range("E13").Value = "  -1.50%  "

$ws.Range("D14").Value = "'5.495"
$ws.Range("E14").Value = "  -0.06%  "

$ws.Range("D15").Value = "'91.93"
$ws.Range("E15").Value = "  -0.03%  "

$ws.Range("D16").Value = "'1.0000"
$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("D17").Value = "'0.000008733"
$ws.Range("E17").Value = "  -1.46%  "

$ws.Range("D18").Value = "'0.9995"
$ws.Range("E18").Value = "  -0.38%  "

$ws.Range("D19").Value = "27.887.85"
$ws.Range("E19").Value = "  -0.60%  "

$ws.Range("E20").Value = "  -1.18%  "

$ws.Range("D21").Value = "'5.181"
$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("E23").Value = "  -0.96%  "

$ws.Range("D24").Value = "'153.04"
$ws.Range("E24").Value = "  -1.60%  "

$ws.Range("D25").Value = "'1.861"
$ws.Range("E25").Value = "  -2.88%  "

$ws.Range("D26").Value = "'2.222"
$ws.Range("E26").Value = "  +3.89%  "

$ws.Range("E27").Value = "  -0.60%  "

$ws.Range("D28").Value = "'115.29"
$ws.Range("E28").Value = "  -1.43%  "

$ws.Range("D29").Value = "'4.907"
$ws.Range("E29").Value = "  -1.24%  "

$ws.Range("D30").Value = "'0.09016"
$ws.Range("E30").Value = "  +0.73%  "

$ws.Range("D31").Value = "'3.197"
$ws.Range("E31").Value = "  -3.12%  "

$ws.Range("D32").Value = "'4.825"
$ws.Range("E32").Value = "  +3.08%  "

$ws.Range("D33").Value = "'1.231"
$ws.Range("E33").Value = "  -2.65%  "

$ws.Range("D34").Value = "'0.7751"
$ws.Range("E34").Value = "  -0.51%  "

$ws.Range("D35").Value = "'0.02089"
$ws.Range("E35").Value = "  +1.65%  "

$ws.Range("D36").Value = "'2.546"
$ws.Range("E36").Value = "  -2.88%  "

$ws.Range("E37").Value = "  -1.59%  "

$ws.Range("E40").Value = "  -0.92%  "

$ws.Range("D41").Value = "'6.902"
$ws.Range("E41").Value = "  -1.79%  "

$ws.Range("D42").Value = "'8.488"
$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("E43").Value = "  -0.89%  "

$ws.Range("D44").Value = "'112.18"
$ws.Range("E44").Value = "  +3.74%  "

$ws.Range("D45").Value = "'10.62"
$ws.Range("E45").Value = "  -0.73%  "

$ws.Range("D46").Value = "'0.4840"
$ws.Range("E46").Value = "  +0.16%  "

$ws.Range("D47").Value = "'0.9994"
$ws.Range("E47").Value = "  -0.46%  "

$ws.Range("D48").Value = "'1.635"
$ws.Range("E48").Value = "  -1.00%  "

$ws.Range("D49").Value = "'67.56"
$ws.Range("E49").Value = "  -0.56%  "

$ws.Range("E50").Value = "  -0.29%  "

$ws.Range("D51").Value = "'0.9080"
$ws.Range("E51").Value = "  +0.86%  "

# Row 38/39 swap: TheSandbox/MXToken positions swap with updated prices
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'3.048"
$ws.Range("E38").Value = "  +1.50%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.5565"
$ws.Range("E39").Value = "  +1.28%  "
